$d = $word.ActiveDocument

# --- Locate "6 (seis) " inside the "jornada laboral" paragraph -------------
$target = $d.Content
$found = $target.Find.Execute("6 (seis) ", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '6 (seis) ' text to replace"
}
$oldStart = $target.Start

# --- Replace the wording: "6 (seis) " -> "cuatro " --------------------------
# (done first, while the run is still a single contiguous range; any later
#  bookmark insertion will cleanly split the surrounding runs without
#  re-merging them)
$target.Text = "cuatro "
$newEnd = $target.End

# --- Split the run boundary right before "cuatro " --------------------------
# Adding + immediately removing a bookmark at a point inside a run forces
# Word to split that run there, without otherwise touching its content.
$boundaryBefore = $d.Range($oldStart, $oldStart)
$d.Bookmarks.Add("_TmpSplit", $boundaryBefore) | Out-Null
$d.Bookmarks("_TmpSplit").Delete()

# --- Re-create the "_GoBack" bookmark right after "cuatro " -----------------
# Word only ever keeps a single "_GoBack" bookmark, tracking the last edited
# spot; adding it here moves it from wherever it used to be (it will be
# removed from its old location automatically since bookmark names are
# unique) and this also splits "cuatro " from "horas de trabajo. ".
$goBackPos = $d.Range($newEnd, $newEnd)
$d.Bookmarks.Add("_GoBack", $goBackPos) | Out-Null
